# Add columns I0 (I) and IF (J) to the worksheet, per the commit "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style used by the other header cells (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-67: I = I0 value, J = IF value ---
$data = @{
    2  = @(9, 9)
    3  = @(8, 8)
    4  = @(10, 10)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(7, 7)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(9, 9)
    19 = @(8, 8)
    20 = @(9, 9)
    21 = @(9, 9)
    22 = @(9, 9)
    23 = @(10, 11)
    24 = @(9, 9)
    25 = @(9, 9)
    26 = @(9, 9)
    27 = @(9, 9)
    28 = @(9, 9)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(9, 9)
    32 = @(9, 9)
    33 = @(9, 9)
    34 = @(9, 9)
    35 = @(9, 9)
    36 = @(9, 9)
    37 = @(9, 9)
    38 = @(9, 9)
    39 = @(10, 10)
    40 = @(9, 9)
    41 = @(8, 8)
    42 = @(9, 9)
    43 = @(9, 9)
    44 = @(10, 10)
    45 = @(9, 9)
    46 = @(9, 9)
    47 = @(10, 10)
    48 = @(9, 9)
    49 = @(8, 8)
    50 = @(9, 9)
    51 = @(9, 9)
    52 = @(9, 9)
    53 = @(9, 9)
    54 = @(9, 9)
    55 = @(9, 9)
    56 = @(9, 9)
    57 = @(9, 9)
    58 = @(10, 10)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(9, 9)
    64 = @(6, 6)
    65 = @(4, 4)
    66 = @(6, 6)
    67 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

$wb.Save()
